$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status column (C) for tracked activities
$ws.Range("C11").Value = "HACIENDO"
$ws.Range("C13").Value = "HACIENDO"
$ws.Range("C14").Value = "HACIENDO"
$ws.Range("C15").Value = "HACIENDO"
$ws.Range("C16").Value = "HACIENDO"
$ws.Range("C17").Value = "HACIENDO"
$ws.Range("C18").Value = "HACER"

# Match the saved selection position from the source edit
$ws.Range("C11").Select()
